# Add a "RNN 모델 이용 " textbox to slide 5 (TextBox 7, shape id 8),
# matching the identical textbox already present on slide 4, and place
# it in the shape tree right before the "Picture 2" picture shape.

$p = $ppt.ActivePresentation
$s4 = $p.Slides.Item(4)
$s5 = $p.Slides.Item(5)

# The slide's shape-id counter is monotonic and never reuses numbers,
# even across add/delete. Slide 5 already has a shape with id=6 (the
# title placeholder), so the counter would otherwise hand out id=7 to
# the next new shape. Burn through throwaway ids 2,3,4,5,7 first so the
# real textbox lands on id=8 - same id the author's copy ended up with
# on this slide (matching slide 4's "TextBox 7").
$burn = @()
for ($i = 1; $i -le 5; $i++) {
    $burn += $s5.Shapes.AddTextbox(1, 0, 0, 10, 10)
}
foreach ($d in $burn) {
    $d.Delete()
}

# Copy the existing "TextBox 7" shape from slide 4 (same text, position,
# size and run-level formatting we need) and paste it onto slide 5.
$src = $s4.Shapes.Item("TextBox 7")
$src.Copy()
$new = $s5.Shapes.Paste().Item(1)

# Paste appends the shape at the end of the z-order; move it so it sits
# right after the title placeholder and right before "Picture 2" (i.e.
# into slot 2 of 3), matching the target shape order.
$new.ZOrder(1)
$new.ZOrder(2)
